# DOMA-2542 Localization for Excel template (ticket_report_status_executor)
#
# The template's third data row uses shared-string placeholders of the form
# "{d.tickets[i + 1].<field>}". This normalises them to "{d.tickets[i+1].<field>}"
# (no spaces around the "+") to match the templating engine's expected syntax,
# matching row 2's "{d.tickets[i].<field>}" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @(
    "categoryClassifier",
    "address",
    "processing",
    "completed",
    "canceled",
    "deferred",
    "closed",
    "new_or_reopened"
)

for ($i = 0; $i -lt $fields.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(3, $col).Value = "{d.tickets[i+1]." + $fields[$i] + "}"
}
